# Apply updated crypto price/volume data to Sheet1 (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.166.07"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.784.86"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.26"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.12"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.791.33"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.98"
$ws.Range("E14").Value = "  -4.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.145.29"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.623"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.68"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.75"
$ws.Range("E19").Value = "  +0.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("E20").Value = "  +1.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.29"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("E27").Value = "  +0.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("E32").Value = "  +1.76%  "

$ws.Range("E33").Value = "  +3.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.448.62"
$ws.Range("E35").Value = "  +3.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.55"
$ws.Range("E36").Value = "  +9.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.654"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("E38").Value = "  +1.57%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.57"
$ws.Range("E40").Value = "  +4.65%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.916"
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.70"
$ws.Range("E44").Value = "  +2.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0518"

$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("E49").Value = "  -5.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.97"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("E51").Value = "  +0.26%  "
